# Applies the scrape update described by the commit
#   "Horarios actualizados Linea 141 - 76"
#
# The workbook has 3 sheets: LP1912, LP1912-215, 6203-6173.
# LP1912 and 6203-6173 receive identical row-level updates (both mirror the
# same underlying scrape; only the title text in A1 differs between them).
# LP1912-215 is a filtered view (only "215*" lines) that only gets two new
# rows appended, plus the header refresh.

$wb = $excel.ActiveWorkbook

function Set-Row5 {
    param($ws, $row, $a, $b, $c, $d, $e)
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

function Swap-Rows4 {
    # Swaps columns A-D between two rows (column E is blank on both sides
    # for every pair that needs swapping in this update).
    param($ws, $row1, $row2)

    $a1 = $ws.Cells.Item($row1, 1).Value2
    $b1 = $ws.Cells.Item($row1, 2).Value2
    $c1 = $ws.Cells.Item($row1, 3).Value2
    $d1 = $ws.Cells.Item($row1, 4).Value2

    $a2 = $ws.Cells.Item($row2, 1).Value2
    $b2 = $ws.Cells.Item($row2, 2).Value2
    $c2 = $ws.Cells.Item($row2, 3).Value2
    $d2 = $ws.Cells.Item($row2, 4).Value2

    $ws.Cells.Item($row1, 1).Value = $a2
    $ws.Cells.Item($row1, 2).Value = $b2
    $ws.Cells.Item($row1, 3).Value = $c2
    $ws.Cells.Item($row1, 4).Value = $d2

    $ws.Cells.Item($row2, 1).Value = $a1
    $ws.Cells.Item($row2, 2).Value = $b1
    $ws.Cells.Item($row2, 3).Value = $c1
    $ws.Cells.Item($row2, 4).Value = $d1
}

function Update-MainSheet {
    # Used for LP1912 and 6203-6173 (245 data rows -> 249 data rows).
    param($ws)

    # --- header block ---
    $ws.Range("A2").Value = "Última actualización: 12:46:01"
    $ws.Range("A3").Value = "Total filas: 249"

    # --- insert the 4 newly scraped rows ---
    # Processed from the bottom of the sheet upward so that every row
    # number referenced below is still the one seen in the original
    # (245-row) layout at the moment it is used - inserting a row only
    # shifts the rows at/after the insertion point, never the ones above it.

    # New last row -> final row 254 (appended after old row 250).
    $ws.Rows.Item(251).Insert()
    Set-Row5 $ws 251 "12:46:01" "14:45" "215B_EL PATO" 119 "LP1912"

    # New row -> final row 251 (inserted right before old row 249).
    $ws.Rows.Item(249).Insert()
    Set-Row5 $ws 249 "12:46:01" "14:34" "215C_LA PLATA" 108 "LP1912"

    # New row -> final row 241 (inserted right before old row 240).
    $ws.Rows.Item(240).Insert()
    Set-Row5 $ws 240 "12:46:01" "13:43" "14_ABASTO" 57 "LP1912"

    # New row -> final row 237 (inserted right before old row 237).
    $ws.Rows.Item(237).Insert()
    Set-Row5 $ws 237 "12:46:01" "13:38" "23_HERNANDEZ" 52 "LP1912"

    # --- re-sorted ties (same Hora_Llegada, new relative scrape order) ---
    Swap-Rows4 $ws 191 192
    Swap-Rows4 $ws 175 176
    Swap-Rows4 $ws 114 115
    Swap-Rows4 $ws 67 68
}

function Update-215Sheet {
    # Used for LP1912-215 (65 data rows -> 67 data rows): only two brand
    # new "215*" rows are appended, matching the two newest 215-rows added
    # to the main sheets above.
    param($ws)

    $ws.Range("A2").Value = "Última actualización: 12:46:01"
    $ws.Range("A3").Value = "Total filas: 67"

    Set-Row5 $ws 71 "12:46:01" "14:34" "215C_LA PLATA" 108 "LP1912"
    Set-Row5 $ws 72 "12:46:01" "14:45" "215B_EL PATO" 119 "LP1912"
}

Update-MainSheet $wb.Worksheets.Item("LP1912")
Update-MainSheet $wb.Worksheets.Item("6203-6173")
Update-215Sheet $wb.Worksheets.Item("LP1912-215")
